$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.710.11'
$ws.Range("E2").Value = '  -2.69%  '
$ws.Range("D3").Value = '1.559.16'
$ws.Range("E3").Value = '  -0.30%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '205.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.486'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.63%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.89'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -0.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0583'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0864'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = '1.781.50'
$ws.Range("E12").Value = '  -0.34%  '
$ws.Range("D13").Value = '1.546.93'
$ws.Range("E13").Value = '  -1.12%  '
$ws.Range("E14").Value = '  -2.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.510'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.22%  '
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '61.48'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.96%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '26.750.44'
$ws.Range("E17").Value = '  -2.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '214.55'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.55%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.32'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.05%  '
$ws.Range("D20").Value = '0.0₃0676'
$ws.Range("E20").Value = '  -1.74%  '
$ws.Range("E21").Value = '  +0.13%  '
$ws.Range("E22").Value = '  -0.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.32'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.23%  '
$ws.Range("E24").Value = '  -1.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.41'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.37%  '
$ws.Range("E26").Value = '  +0.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.86'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.95%  '
$ws.Range("E28").Value = '  +0.15%  '
$ws.Range("E29").Value = '  -1.47%  '
$ws.Range("E30").Value = '  -3.49%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0461'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.76%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.15'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.50%  '
$ws.Range("D33").Value = '1.385.48'
$ws.Range("E33").Value = '  +1.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.90'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.44%  '
$ws.Range("E35").Value = '  +1.64%  '
$ws.Range("E36").Value = '  -0.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.931'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.32%  '
$ws.Range("E38").Value = '  -2.89%  '
$ws.Range("E39").Value = '  -1.77%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.514'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.28%  '
$ws.Range("E41").Value = '  +0.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.997'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.32%  '
$ws.Range("E43").Value = '  +2.64%  '
$ws.Range("E44").Value = '  +1.82%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.20'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.40%  '
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.75'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.13%  '
$ws.Range("D47").Value = '1.694.63'
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.58'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.25%  '
$ws.Range("D49").Value = '0.0₇0970'
$ws.Range("E49").Value = '  -1.67%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0492'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0944'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.89%  '
